# Update the dated worksheet heading and refresh the 25 practice problems
# (3-digit x 1-digit multiplication) with the new day's values/answers.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-15 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-16 Thursday", 2) | Out-Null
$d.Content.Find.Execute("685×2=1370", $true, $false, $false, $false, $false, $true, 1, $false, "486×5=2430", 2) | Out-Null
$d.Content.Find.Execute("481×4=1924", $true, $false, $false, $false, $false, $true, 1, $false, "837×7=5859", 2) | Out-Null
$d.Content.Find.Execute("951×8=7608", $true, $false, $false, $false, $false, $true, 1, $false, "486×9=4374", 2) | Out-Null
$d.Content.Find.Execute("608×3=1824", $true, $false, $false, $false, $false, $true, 1, $false, "750×4=3000", 2) | Out-Null
$d.Content.Find.Execute("329×5=1645", $true, $false, $false, $false, $false, $true, 1, $false, "981×8=7848", 2) | Out-Null
$d.Content.Find.Execute("973×2=1946", $true, $false, $false, $false, $false, $true, 1, $false, "478×3=1434", 2) | Out-Null
$d.Content.Find.Execute("150×9=1350", $true, $false, $false, $false, $false, $true, 1, $false, "609×7=4263", 2) | Out-Null
$d.Content.Find.Execute("663×8=5304", $true, $false, $false, $false, $false, $true, 1, $false, "702×2=1404", 2) | Out-Null
$d.Content.Find.Execute("827×8=6616", $true, $false, $false, $false, $false, $true, 1, $false, "577×5=2885", 2) | Out-Null
$d.Content.Find.Execute("255×3=765", $true, $false, $false, $false, $false, $true, 1, $false, "719×8=5752", 2) | Out-Null
$d.Content.Find.Execute("264×3=792", $true, $false, $false, $false, $false, $true, 1, $false, "612×3=1836", 2) | Out-Null
$d.Content.Find.Execute("641×6=3846", $true, $false, $false, $false, $false, $true, 1, $false, "137×9=1233", 2) | Out-Null
$d.Content.Find.Execute("891×3=2673", $true, $false, $false, $false, $false, $true, 1, $false, "141×9=1269", 2) | Out-Null
$d.Content.Find.Execute("169×7=1183", $true, $false, $false, $false, $false, $true, 1, $false, "371×2=742", 2) | Out-Null
$d.Content.Find.Execute("348×4=1392", $true, $false, $false, $false, $false, $true, 1, $false, "308×9=2772", 2) | Out-Null
$d.Content.Find.Execute("377×5=1885", $true, $false, $false, $false, $false, $true, 1, $false, "714×8=5712", 2) | Out-Null
$d.Content.Find.Execute("840×5=4200", $true, $false, $false, $false, $false, $true, 1, $false, "971×6=5826", 2) | Out-Null
$d.Content.Find.Execute("956×9=8604", $true, $false, $false, $false, $false, $true, 1, $false, "256×5=1280", 2) | Out-Null
$d.Content.Find.Execute("994×2=1988", $true, $false, $false, $false, $false, $true, 1, $false, "640×5=3200", 2) | Out-Null
$d.Content.Find.Execute("712×7=4984", $true, $false, $false, $false, $false, $true, 1, $false, "206×4=824", 2) | Out-Null
$d.Content.Find.Execute("671×6=4026", $true, $false, $false, $false, $false, $true, 1, $false, "416×2=832", 2) | Out-Null
$d.Content.Find.Execute("276×5=1380", $true, $false, $false, $false, $false, $true, 1, $false, "720×5=3600", 2) | Out-Null
$d.Content.Find.Execute("951×2=1902", $true, $false, $false, $false, $false, $true, 1, $false, "399×4=1596", 2) | Out-Null
$d.Content.Find.Execute("150×4=600", $true, $false, $false, $false, $false, $true, 1, $false, "830×6=4980", 2) | Out-Null
$d.Content.Find.Execute("473×2=946", $true, $false, $false, $false, $false, $true, 1, $false, "730×3=2190", 2) | Out-Null
